# "Move data to index.js" -- reposition/resize the Dialogs/DialogsList
# shape and its connected chain of shapes+connectors on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$EMU_PER_POINT = 12700

# PowerPoint's Shape.Left/Top/Width/Height are single-precision (float32)
# point values. Converting an EMU value to points and back truncates
# towards zero, so a plain "$emu / 12700" can land 1 EMU short of the
# intended value after the float32 round-trip. Biasing by half an EMU
# keeps the float32 value safely inside the correct EMU bucket.
function ToPt($emu) {
    if ($emu -eq 0) { return 0.0 }
    return ($emu + 0.5) / $EMU_PER_POINT
}

# Shape 97: "Dialogs" -> "DialogsList", widened/moved left.
$sh97 = Get-ShapeById $s 97
$sh97.Left = ToPt 3169025
$sh97.Top = ToPt 2447125
$sh97.Width = ToPt 815400
$sh97.Height = ToPt 248100
$sh97.TextFrame.TextRange.Text = "DialogsList"

# Shape 99: "DialogItem" moved left (size unchanged).
$sh99 = Get-ShapeById $s 99
$sh99.Left = ToPt 3169025
$sh99.Top = ToPt 3285475
$sh99.Width = ToPt 815400
$sh99.Height = ToPt 248100

# Connector 100 (97 idx2 -> 99 idx0): no longer flipped, now vertical.
$sh100 = Get-ShapeById $s 100
$sh100.HorizontalFlip = 0
$sh100.Left = ToPt 3576725
$sh100.Top = ToPt 2695225
$sh100.Width = ToPt 0
$sh100.Height = ToPt 590400

# Shape 101: "MessageItem" moved left (size unchanged).
$sh101 = Get-ShapeById $s 101
$sh101.Left = ToPt 4161425
$sh101.Top = ToPt 3285475
$sh101.Width = ToPt 947700
$sh101.Height = ToPt 248100

# Connector 102 (98 idx2 -> 101 idx0): now vertical (width -> 0).
$sh102 = Get-ShapeById $s 102
$sh102.Left = ToPt 4635275
$sh102.Top = ToPt 2695225
$sh102.Width = ToPt 0
$sh102.Height = ToPt 590400

# Connector 103 (71 idx2 -> 97 idx0): stays flipped, moved/widened.
$sh103 = Get-ShapeById $s 103
$sh103.Left = ToPt 3576725
$sh103.Top = ToPt 2133425
$sh103.Width = ToPt 537900
$sh103.Height = ToPt 313800
